$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.1172255
$ws.Range("H2").Value = 4.234451
$ws.Range("I2").Value = 0.006295392006363395
$ws.Range("J2").Value = 0.004213039461358209
$ws.Range("M2").Value = 0.06913800000000001
$ws.Range("N2").Value = 0.138276
$ws.Range("O2").Value = 0.0001995195529840091
$ws.Range("P2").Value = 0.000133021882144827
$ws.Range("Q2").Value = 0.146380736619
$ws.Range("R2").Value = 0.585522946476
$ws.Range("S2").Value = 0.000001256053798968729
$ws.Range("T2").Value = 0.0000005604264387002968
$ws.Range("G3").Value = 2.1172255
$ws.Range("H3").Value = 4.234451
$ws.Range("I3").Value = 0.006295392006363395
$ws.Range("J3").Value = 0.004213039461358209
$ws.Range("M3").Value = 76.57257800000001
$ws.Range("N3").Value = 229.717734
$ws.Range("O3").Value = 0.2209743778152849
$ws.Range("P3").Value = 0.2209890750291063
$ws.Range("Q3").Value = 162.121414742339
$ws.Range("R3").Value = 972.728488454034
$ws.Range("S3").Value = 0.00139112033170947
$ws.Range("T3").Value = 0.000931035693626675
$ws.Range("G4").Value = 2.1172255
$ws.Range("H4").Value = 4.234451
$ws.Range("I4").Value = 0.006295392006363395
$ws.Range("J4").Value = 0.004213039461358209
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 267.3435366666667
$ws.Range("N4").Value = 802.03061
$ws.Range("O4").Value = 0.7715042802640716
$ws.Range("P4").Value = 0.7715555937397933
$ws.Range("Q4").Value = 566.0265530908517
$ws.Range("R4").Value = 3396.15931854511
$ws.Range("S4").Value = 0.004856921878849581
$ws.Range("T4").Value = 0.003250594163057412
$ws.Range("G5").Value = 2.1172255
$ws.Range("H5").Value = 4.234451
$ws.Range("I5").Value = 0.006295392006363395
$ws.Range("J5").Value = 0.004213039461358209
$ws.Range("M5").Value = 2.537175666666667
$ws.Range("N5").Value = 7.611527
$ws.Range("O5").Value = 0.007321822367659443
$ws.Range("P5").Value = 0.007322309348955482
$ws.Range("Q5").Value = 5.371773019446167
$ws.Range("R5").Value = 32.230638116677
$ws.Range("S5").Value = 0.00004609374200537597
$ws.Range("T5").Value = 0.00003084917823542158
$ws.Range("I6").Value = 0.004671287948366863
$ws.Range("J6").Value = 0.004689220411201648
$ws.Range("M6").Value = 0.06913800000000001
$ws.Range("N6").Value = 0.138276
$ws.Range("O6").Value = 0.0001995195529840091
$ws.Range("P6").Value = 0.000133021882144827
$ws.Range("Q6").Value = 0.108616996392
$ws.Range("R6").Value = 0.6517019783520001
$ws.Range("S6").Value = 0.0000009320132833177454
$ws.Range("T6").Value = 0.0000006237689248899826
$ws.Range("I7").Value = 0.004671287948366863
$ws.Range("J7").Value = 0.004689220411201648
$ws.Range("M7").Value = 76.57257800000001
$ws.Range("N7").Value = 229.717734
$ws.Range("O7").Value = 0.2209743778152849
$ws.Range("P7").Value = 0.2209890750291063
$ws.Range("Q7").Value = 120.2968472960187
$ws.Range("R7").Value = 1082.671625664168
$ws.Range("S7").Value = 0.001032234947986406
$ws.Range("T7").Value = 0.001036266481279058
$ws.Range("I8").Value = 0.004671287948366863
$ws.Range("J8").Value = 0.004689220411201648
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 267.3435366666667
$ws.Range("N8").Value = 802.03061
$ws.Range("O8").Value = 0.7715042802640716
$ws.Range("P8").Value = 0.7715555937397933
$ws.Range("Q8").Value = 420.001330057969
$ws.Range("R8").Value = 3780.01197052172
$ws.Range("S8").Value = 0.003603918646511008
$ws.Range("T8").Value = 0.003617994238541446
$ws.Range("I9").Value = 0.004671287948366863
$ws.Range("J9").Value = 0.004689220411201648
$ws.Range("M9").Value = 2.537175666666667
$ws.Range("N9").Value = 7.611527
$ws.Range("O9").Value = 0.007321822367659443
$ws.Range("P9").Value = 0.007322309348955482
$ws.Range("Q9").Value = 3.985946950044889
$ws.Range("R9").Value = 35.873522550404
$ws.Range("S9").Value = 0.00003420234058613049
$ws.Range("T9").Value = 0.00003433592245625469
$ws.Range("G10").Value = 75.40439600000001
$ws.Range("H10").Value = 226.213188
$ws.Range("I10").Value = 0.2242086314485916
$ws.Range("J10").Value = 0.2250693390296979
$ws.Range("M10").Value = 0.06913800000000001
$ws.Range("N10").Value = 0.138276
$ws.Range("O10").Value = 0.0001995195529840091
$ws.Range("P10").Value = 0.000133021882144827
$ws.Range("Q10").Value = 5.213309130648001
$ws.Range("R10").Value = 31.279854783888
$ws.Range("S10").Value = 0.00004473400592177943
$ws.Range("T10").Value = 0.00002993914709082257
$ws.Range("G11").Value = 75.40439600000001
$ws.Range("H11").Value = 226.213188
$ws.Range("I11").Value = 0.2242086314485916
$ws.Range("J11").Value = 0.2250693390296979
$ws.Range("M11").Value = 76.57257800000001
$ws.Range("N11").Value = 229.717734
$ws.Range("O11").Value = 0.2209743778152849
$ws.Range("P11").Value = 0.2209890750291063
$ws.Range("Q11").Value = 5773.908994252889
$ws.Range("R11").Value = 51965.18094827599
$ws.Range("S11").Value = 0.04954436283516905
$ws.Range("T11").Value = 0.04973786504958528
$ws.Range("G12").Value = 75.40439600000001
$ws.Range("H12").Value = 226.213188
$ws.Range("I12").Value = 0.2242086314485916
$ws.Range("J12").Value = 0.2250693390296979
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 267.3435366666667
$ws.Range("N12").Value = 802.03061
$ws.Range("O12").Value = 0.7715042802640716
$ws.Range("P12").Value = 0.7715555937397933
$ws.Range("Q12").Value = 20158.87790685386
$ws.Range("R12").Value = 181429.9011616847
$ws.Range("S12").Value = 0.1729779188347381
$ws.Range("T12").Value = 0.1736535075076814
$ws.Range("G13").Value = 75.40439600000001
$ws.Range("H13").Value = 226.213188
$ws.Range("I13").Value = 0.2242086314485916
$ws.Range("J13").Value = 0.2250693390296979
$ws.Range("M13").Value = 2.537175666666667
$ws.Range("N13").Value = 7.611527
$ws.Range("O13").Value = 0.007321822367659443
$ws.Range("P13").Value = 0.007322309348955482
$ws.Range("Q13").Value = 191.3141986908973
$ws.Range("R13").Value = 1721.827788218076
$ws.Range("S13").Value = 0.00164161577276261
$ws.Range("T13").Value = 0.001648027325340388
$ws.Range("G14").Value = 1.7411535
$ws.Range("H14").Value = 3.482307
$ws.Range("I14").Value = 0.005177173534775417
$ws.Range("J14").Value = 0.003464698684094803
$ws.Range("M14").Value = 0.06913800000000001
$ws.Range("N14").Value = 0.138276
$ws.Range("O14").Value = 0.0001995195529840091
$ws.Range("P14").Value = 0.000133021882144827
$ws.Range("Q14").Value = 0.120379870683
$ws.Range("R14").Value = 0.4815194827320001
$ws.Range("S14").Value = 0.000001032947349379033
$ws.Range("T14").Value = 0.0000004608807400229959
$ws.Range("G15").Value = 1.7411535
$ws.Range("H15").Value = 3.482307
$ws.Range("I15").Value = 0.005177173534775417
$ws.Range("J15").Value = 0.003464698684094803
$ws.Range("M15").Value = 76.57257800000001
$ws.Range("N15").Value = 229.717734
$ws.Range("O15").Value = 0.2209743778152849
$ws.Range("P15").Value = 0.2209890750291063
$ws.Range("Q15").Value = 133.324612188723
$ws.Range("R15").Value = 799.9476731323381
$ws.Range("S15").Value = 0.001144022700688757
$ws.Range("T15").Value = 0.0007656605574526724
$ws.Range("G16").Value = 1.7411535
$ws.Range("H16").Value = 3.482307
$ws.Range("I16").Value = 0.005177173534775417
$ws.Range("J16").Value = 0.003464698684094803
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 267.3435366666667
$ws.Range("N16").Value = 802.03061
$ws.Range("O16").Value = 0.7715042802640716
$ws.Range("P16").Value = 0.7715555937397933
$ws.Range("Q16").Value = 465.4861345695451
$ws.Range("R16").Value = 2792.91680741727
$ws.Range("S16").Value = 0.003994211541749107
$ws.Range("T16").Value = 0.002673207650336246
$ws.Range("G17").Value = 1.7411535
$ws.Range("H17").Value = 3.482307
$ws.Range("I17").Value = 0.005177173534775417
$ws.Range("J17").Value = 0.003464698684094803
$ws.Range("M17").Value = 2.537175666666667
$ws.Range("N17").Value = 7.611527
$ws.Range("O17").Value = 0.007321822367659443
$ws.Range("P17").Value = 0.007322309348955482
$ws.Range("Q17").Value = 4.4176122921315
$ws.Range("R17").Value = 26.505673752789
$ws.Range("S17").Value = 0.00003790634498817315
$ws.Range("T17").Value = 0.00002536959556586113
$ws.Range("G18").Value = 197.2895866666667
$ws.Range("H18").Value = 591.86876
$ws.Range("I18").Value = 0.586623997698909
$ws.Range("J18").Value = 0.5888759704209946
$ws.Range("M18").Value = 0.06913800000000001
$ws.Range("N18").Value = 0.138276
$ws.Range("O18").Value = 0.0001995195529840091
$ws.Range("P18").Value = 0.000133021882144827
$ws.Range("Q18").Value = 13.64020744296
$ws.Range("R18").Value = 81.84124465776
$ws.Range("S18").Value = 0.0001170429577905787
$ws.Range("T18").Value = 0.00007833338993526215
$ws.Range("G19").Value = 197.2895866666667
$ws.Range("H19").Value = 591.86876
$ws.Range("I19").Value = 0.586623997698909
$ws.Range("J19").Value = 0.5888759704209946
$ws.Range("M19").Value = 76.57257800000001
$ws.Range("N19").Value = 229.717734
$ws.Range("O19").Value = 0.2209743778152849
$ws.Range("P19").Value = 0.2209890750291063
$ws.Range("Q19").Value = 15106.97226362109
$ws.Range("R19").Value = 135962.7503725898
$ws.Range("S19").Value = 0.1296288729030316
$ws.Range("T19").Value = 0.130135156010203
$ws.Range("G20").Value = 197.2895866666667
$ws.Range("H20").Value = 591.86876
$ws.Range("I20").Value = 0.586623997698909
$ws.Range("J20").Value = 0.5888759704209946
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 267.3435366666667
$ws.Range("N20").Value = 802.03061
$ws.Range("O20").Value = 0.7715042802640716
$ws.Range("P20").Value = 0.7715555937397933
$ws.Range("Q20").Value = 52744.09584697151
$ws.Range("R20").Value = 474696.8626227436
$ws.Range("S20").Value = 0.4525829251303292
$ws.Range("T20").Value = 0.4543505489972675
$ws.Range("G21").Value = 197.2895866666667
$ws.Range("H21").Value = 591.86876
$ws.Range("I21").Value = 0.586623997698909
$ws.Range("J21").Value = 0.5888759704209946
$ws.Range("M21").Value = 2.537175666666667
$ws.Range("N21").Value = 7.611527
$ws.Range("O21").Value = 0.007321822367659443
$ws.Range("P21").Value = 0.007322309348955482
$ws.Range("Q21").Value = 500.5583385773911
$ws.Range("R21").Value = 4505.025047196519
$ws.Range("S21").Value = 0.004295156707757674
$ws.Range("T21").Value = 0.00431193202358888
$ws.Range("G22").Value = 58.19014966666668
$ws.Range("H22").Value = 174.570449
$ws.Range("I22").Value = 0.1730235173629937
$ws.Range("J22").Value = 0.1736877319926528
$ws.Range("M22").Value = 0.06913800000000001
$ws.Range("N22").Value = 0.138276
$ws.Range("O22").Value = 0.0001995195529840091
$ws.Range("P22").Value = 0.000133021882144827
$ws.Range("Q22").Value = 4.023150567654001
$ws.Range("R22").Value = 24.138903405924
$ws.Range("S22").Value = 0.00003452157483998543
$ws.Range("T22").Value = 0.00002310426901512895
$ws.Range("G23").Value = 58.19014966666668
$ws.Range("H23").Value = 174.570449
$ws.Range("I23").Value = 0.1730235173629937
$ws.Range("J23").Value = 0.1736877319926528
$ws.Range("M23").Value = 76.57257800000001
$ws.Range("N23").Value = 229.717734
$ws.Range("O23").Value = 0.2209743778152849
$ws.Range("P23").Value = 0.2209890750291063
$ws.Range("Q23").Value = 4455.769774182509
$ws.Range("R23").Value = 40101.92796764257
$ws.Range("S23").Value = 0.03823376409669967
$ws.Range("T23").Value = 0.03838309123695967
$ws.Range("G24").Value = 58.19014966666668
$ws.Range("H24").Value = 174.570449
$ws.Range("I24").Value = 0.1730235173629937
$ws.Range("J24").Value = 0.1736877319926528
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 267.3435366666667
$ws.Range("N24").Value = 802.03061
$ws.Range("O24").Value = 0.7715042802640716
$ws.Range("P24").Value = 0.7715555937397933
$ws.Range("Q24").Value = 15556.76041104932
$ws.Range("R24").Value = 140010.8436994439
$ws.Range("S24").Value = 0.1334883842318945
$ws.Range("T24").Value = 0.1340097411829093
$ws.Range("G25").Value = 58.19014966666668
$ws.Range("H25").Value = 174.570449
$ws.Range("I25").Value = 0.1730235173629937
$ws.Range("J25").Value = 0.1736877319926528
$ws.Range("M25").Value = 2.537175666666667
$ws.Range("N25").Value = 7.611527
$ws.Range("O25").Value = 0.007321822367659443
$ws.Range("P25").Value = 0.007322309348955482
$ws.Range("Q25").Value = 147.6386317739581
$ws.Range("R25").Value = 1328.747685965623
$ws.Range("S25").Value = 0.001266847459559479
$ws.Range("T25").Value = 0.001271795303768676
